$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 138
$ws.Range("I9").Value = 133
$ws.Range("K9").Value = 133
$ws.Range("M9").Value = 36

$ws.Range("H19").Value = 1166.75
$ws.Range("I19").Value = 1315.3334
$ws.Range("J19").Value = 721
$ws.Range("K19").Value = 1315.3334
$ws.Range("L19").Value = 721
$ws.Range("M19").Value = -1140.3334
$ws.Range("N19").Value = -1071

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H51").Value = 2984.4722
$ws.Range("I51").Value = 2975.3489
$ws.Range("J51").Value = 2998
$ws.Range("K51").Value = 2975.3489
$ws.Range("L51").Value = 2998
$ws.Range("M51").Value = -2491.3489
$ws.Range("N51").Value = -3966

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H112").Value = 2142.2104
$ws.Range("J112").Value = 2842.3333
$ws.Range("L112").Value = 8526.999899999999
$ws.Range("N112").Value = -10742.9999

$ws.Range("H116").Value = 20772.467
$ws.Range("I116").Value = 4566.1665
$ws.Range("J116").Value = 31576.666
$ws.Range("K116").Value = 4566.1665
$ws.Range("L116").Value = 31576.666
$ws.Range("M116").Value = -1124.1665
$ws.Range("N116").Value = -38460.666

$ws.Range("H132").Value = 11245.4
$ws.Range("I132").Value = 14023.059
$ws.Range("K132").Value = 42069.177
$ws.Range("M132").Value = -39539.177

$ws.Range("H138").Value = 2216.3333
$ws.Range("I138").Value = 1683.2593
$ws.Range("J138").Value = 2498.549
$ws.Range("K138").Value = 5049.7779
$ws.Range("L138").Value = 7495.647
$ws.Range("M138").Value = 90.22209999999995
$ws.Range("N138").Value = -17775.647

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1491.5
$ws.Range("I2").Value = 787.5
$ws.Range("J2").Value = 2899.5
$ws.Range("K2").Value = 787.5
$ws.Range("L2").Value = 2899.5
$ws.Range("M2").Value = -674.5
$ws.Range("N2").Value = -3125.5

$ws.Range("H32").Value = 124426.23
$ws.Range("I32").Value = 133134.83
$ws.Range("K32").Value = 133134.83
$ws.Range("M32").Value = -132847.83

$ws.Range("H45").Value = 5949.5
$ws.Range("I45").Value = 4999
$ws.Range("K45").Value = 4999
$ws.Range("M45").Value = -4622

$ws.Range("H63").Value = 2408.818
$ws.Range("I63").Value = 2408.818
$ws.Range("K63").Value = 2408.818
$ws.Range("M63").Value = -1722.818

$ws.Range("H66").Value = 2408.818
$ws.Range("I66").Value = 2408.818
$ws.Range("K66").Value = 12044.09
$ws.Range("M66").Value = -8612.09

$ws.Range("H76").Value = 29999
$ws.Range("J76").Value = 29999
$ws.Range("L76").Value = 29999
$ws.Range("N76").Value = -30675

$ws.Range("H79").Value = 29999
$ws.Range("J79").Value = 29999
$ws.Range("L79").Value = 29999
$ws.Range("N79").Value = -32339

$ws.Range("H116").Value = 1491.5
$ws.Range("I116").Value = 787.5
$ws.Range("J116").Value = 2899.5
$ws.Range("K116").Value = 787.5
$ws.Range("L116").Value = 2899.5
$ws.Range("M116").Value = 1506.5
$ws.Range("N116").Value = -7487.5

$ws.Range("H132").Value = 2502489.8
$ws.Range("I132").Value = 2780320.5
$ws.Range("J132").Value = 2014
$ws.Range("K132").Value = 8340961.5
$ws.Range("L132").Value = 6042
$ws.Range("M132").Value = -8338431.5
$ws.Range("N132").Value = -11102

$ws.Range("H134").Value = 64997
$ws.Range("J134").Value = 64997
$ws.Range("L134").Value = 64997
$ws.Range("N134").Value = -75137

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1491.5
$ws.Range("I3").Value = 787.5
$ws.Range("J3").Value = 2899.5
$ws.Range("K3").Value = 787.5
$ws.Range("L3").Value = 2899.5
$ws.Range("M3").Value = -673.5
$ws.Range("N3").Value = -3127.5

$ws.Range("H80").Value = 1909.7273
$ws.Range("I80").Value = 1409.3334
$ws.Range("K80").Value = 1409.3334
$ws.Range("M80").Value = -411.3334

$ws.Range("H83").Value = 1909.7273
$ws.Range("I83").Value = 1409.3334
$ws.Range("K83").Value = 7046.666999999999
$ws.Range("M83").Value = -2054.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9155

$ws.Range("H86").Value = 10526.134
$ws.Range("I86").Value = 17642.143
$ws.Range("J86").Value = 4299.625
$ws.Range("K86").Value = 17642.143
$ws.Range("L86").Value = 4299.625
$ws.Range("M86").Value = -16519.143
$ws.Range("N86").Value = -6545.625

$ws.Range("H89").Value = 10526.134
$ws.Range("I89").Value = 17642.143
$ws.Range("J89").Value = 4299.625
$ws.Range("K89").Value = 88210.715
$ws.Range("L89").Value = 21498.125
$ws.Range("M89").Value = -82594.715
$ws.Range("N89").Value = -32730.125

$ws.Range("H116").Value = 74337
$ws.Range("J116").Value = 74000
$ws.Range("L116").Value = 74000
$ws.Range("N116").Value = -83178

$ws.Range("H134").Value = 3022.7144
$ws.Range("I134").Value = 2424.5908
$ws.Range("K134").Value = 7273.7724
$ws.Range("M134").Value = -4738.7724

$ws.Range("H141").Value = 414319.6
$ws.Range("I141").Value = 23000
$ws.Range("K141").Value = 23000
$ws.Range("M141").Value = -17820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 332.57144
$ws.Range("J23").Value = 83.666664
$ws.Range("L23").Value = 250.999992
$ws.Range("N23").Value = -720.999992

$ws.Range("H55").Value = 56004388
$ws.Range("I55").Value = 210000110
$ws.Range("J55").Value = 5943.636
$ws.Range("K55").Value = 630000330
$ws.Range("L55").Value = 17830.908
$ws.Range("M55").Value = -630000153
$ws.Range("N55").Value = -18184.908

$ws.Range("H58").Value = 6999.9
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 3000
$ws.Range("M58").Value = -2872

$ws.Range("H64").Value = 8083.4736
$ws.Range("I64").Value = 4136.125
$ws.Range("K64").Value = 12408.375
$ws.Range("M64").Value = -12138.375

$ws.Range("H67").Value = 8083.4736
$ws.Range("I67").Value = 4136.125
$ws.Range("K67").Value = 12408.375
$ws.Range("M67").Value = -11472.375

$ws.Range("H112").Value = 11327.286
$ws.Range("I112").Value = 4430.5
$ws.Range("K112").Value = 13291.5
$ws.Range("M112").Value = -12183.5

$ws.Range("H123").Value = 12281.6875
$ws.Range("J123").Value = 14549.7
$ws.Range("L123").Value = 43649.10000000001
$ws.Range("N123").Value = -48549.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 9000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H122").Value = 66576.94
$ws.Range("I122").Value = 94566.63
$ws.Range("J122").Value = 4999.6
$ws.Range("K122").Value = 283699.89
$ws.Range("L122").Value = 14998.8
$ws.Range("M122").Value = -281249.89
$ws.Range("N122").Value = -19898.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 899.2941
$ws.Range("I82").Value = 811.5
$ws.Range("K82").Value = 811.5
$ws.Range("M82").Value = -450.5

$ws.Range("H85").Value = 899.2941
$ws.Range("I85").Value = 811.5
$ws.Range("K85").Value = 811.5
$ws.Range("M85").Value = 436.5

$ws.Range("H122").Value = 4294.294
$ws.Range("I122").Value = 3937.6875
$ws.Range("K122").Value = 11813.0625
$ws.Range("M122").Value = -9363.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H136").Value = 3856
$ws.Range("I136").Value = 3713
$ws.Range("K136").Value = 11139
$ws.Range("M136").Value = -8589
